$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card15")

# Add the new "Correction " column by cloning column M (same header
# formatting + the same blank placeholder cells used for the rest of the
# table), then overwrite the header text for the new column.
$ws.Range("M1:M12").Copy($ws.Range("N1:N12"))
$ws.Range("N1").Value = "Correction "

# The existing M1 header ("Event ") loses its trailing space.
$ws.Range("M1").Value = "Event"

# Column M rows 2-12 were blank placeholders; they now read "nan" like the
# rest of the row.
$ws.Range("M2:M12").Value = "nan"
